$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: for each destination row (3-13), the source row in the
# ORIGINAL (pre-edit) sheet whose B,D,E,F,G,H values should land there.
# Row 4 is unchanged (maps to itself); columns A and C never change.
$map = @{
    3  = 6
    4  = 4
    5  = 3
    6  = 5
    7  = 8
    8  = 9
    9  = 11
    10 = 7
    11 = 10
    12 = 13
    13 = 12
}

# Snapshot the original B, D, E, F, G, H values for rows 3-13 before
# writing anything, since several rows read from each other.
$orig = @{}
foreach ($r in 2..13) {
    $orig[$r] = @{
        B = $ws.Cells.Item($r, 2).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $vals = $orig[$srcRow]

    $ws.Cells.Item($destRow, 2).Value = $vals.B

    $ws.Cells.Item($destRow, 4).Value = $vals.D
    $ws.Cells.Item($destRow, 5).Value = $vals.E
    $ws.Cells.Item($destRow, 6).Value = $vals.F

    if ($null -eq $vals.G) {
        $ws.Cells.Item($destRow, 7).Value = ""
    } else {
        $ws.Cells.Item($destRow, 7).Value = $vals.G
    }

    if ($null -eq $vals.H) {
        $ws.Cells.Item($destRow, 8).Value = ""
    } else {
        $ws.Cells.Item($destRow, 8).Value = $vals.H
    }
}
